# Auto-generated update: meteocat daily summary refresh (2026-02-23 14:50 run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-23 14:48:15'
$ws.Range("K2").Value = '11.7 MJ/m2'
$ws.Range("M2").Value = '12.8 °C 14:29 TU'
$ws.Range("O2").Value = '5.9 °C'

$ws.Range("E3").Value = '2026-02-23 14:48:17'
$ws.Range("H3").Value = "'31%"
$ws.Range("K3").Value = '13.7 MJ/m2'
$ws.Range("O3").Value = '3.4 °C'

$ws.Range("E4").Value = '2026-02-23 14:48:19'
$ws.Range("H4").Value = "'69%"
$ws.Range("J4").Value = '1025.3 hPa'
$ws.Range("K4").Value = '12.5 MJ/m2'
$ws.Range("L4").Value = '18.7 km/h - 105º 14:08 TU'
$ws.Range("M4").Value = '21.6 °C 14:08 TU'
$ws.Range("O4").Value = '11.1 °C'

$ws.Range("E5").Value = '2026-02-23 14:48:22'
$ws.Range("I5").Value = '0.2 mm'
$ws.Range("K5").Value = '13.4 MJ/m2'
$ws.Range("M5").Value = '9.0 °C 14:01 TU'
$ws.Range("O5").Value = '3.9 °C'

$ws.Range("E6").Value = '2026-02-23 14:48:24'
$ws.Range("J6").Value = '1025.0 hPa'
$ws.Range("K6").Value = '12.9 MJ/m2'
$ws.Range("O6").Value = '13.1 °C'

$ws.Range("E7").Value = '2026-02-23 14:48:27'
$ws.Range("J7").Value = '1025.1 hPa'
$ws.Range("K7").Value = '13.2 MJ/m2'
$ws.Range("O7").Value = '13.8 °C'

$ws.Range("E8").Value = '2026-02-23 14:48:29'
$ws.Range("K8").Value = '13.2 MJ/m2'
$ws.Range("O8").Value = '14.0 °C'

$ws.Range("E9").Value = '2026-02-23 14:48:32'
$ws.Range("H9").Value = "'74%"
$ws.Range("K9").Value = '12.7 MJ/m2'
$ws.Range("L9").Value = '15.1 km/h - 158º 14:01 TU'
$ws.Range("O9").Value = '11.6 °C'

$ws.Range("E10").Value = '2026-02-23 14:48:34'
$ws.Range("H10").Value = "'79%"
$ws.Range("K10").Value = '13.0 MJ/m2'
$ws.Range("M10").Value = '20.4 °C 14:29 TU'
$ws.Range("O10").Value = '9.9 °C'

$ws.Range("E11").Value = '2026-02-23 14:48:36'
$ws.Range("H11").Value = "'74%"
$ws.Range("O11").Value = '7.4 °C'

$ws.Range("E12").Value = '2026-02-23 14:48:39'
$ws.Range("O12").Value = '9.9 °C'

$ws.Range("E13").Value = '2026-02-23 14:48:41'
$ws.Range("H13").Value = "'68%"
$ws.Range("J13").Value = '1028.8 hPa'
$ws.Range("K13").Value = '13.1 MJ/m2'
$ws.Range("M13").Value = '21.7 °C 14:08 TU'
$ws.Range("O13").Value = '4.5 °C'

$ws.Range("E14").Value = '2026-02-23 14:48:44'
$ws.Range("K14").Value = '12.6 MJ/m2'
$ws.Range("O14").Value = '12.5 °C'

$ws.Range("E15").Value = '2026-02-23 14:48:46'
$ws.Range("H15").Value = "'71%"
$ws.Range("O15").Value = '11.8 °C'

$ws.Range("E16").Value = '2026-02-23 14:48:48'
$ws.Range("H16").Value = "'16%"
$ws.Range("K16").Value = '11.0 MJ/m2'
$ws.Range("M16").Value = '6.8 °C 14:10 TU'
$ws.Range("O16").Value = '3.5 °C'

$ws.Range("E17").Value = '2026-02-23 14:48:51'
$ws.Range("H17").Value = "'45%"
$ws.Range("K17").Value = '14.4 MJ/m2'
$ws.Range("M17").Value = '12.4 °C 14:24 TU'
$ws.Range("O17").Value = '8.5 °C'

$ws.Range("E18").Value = '2026-02-23 14:48:53'
$ws.Range("H18").Value = "'79%"
$ws.Range("J18").Value = '1025.6 hPa'
$ws.Range("K18").Value = '13.2 MJ/m2'
$ws.Range("M18").Value = '21.4 °C 14:28 TU'
$ws.Range("O18").Value = '9.2 °C'

$ws.Range("E19").Value = '2026-02-23 14:48:55'
$ws.Range("K19").Value = '12.8 MJ/m2'
$ws.Range("M19").Value = '17.5 °C 14:19 TU'
$ws.Range("O19").Value = '11.9 °C'

$ws.Range("E20").Value = '2026-02-23 14:48:58'
$ws.Range("K20").Value = '13.8 MJ/m2'
$ws.Range("O20").Value = '4.2 °C'

$ws.Range("E21").Value = '2026-02-23 14:49:00'
$ws.Range("H21").Value = "'64%"
$ws.Range("J21").Value = '1027.2 hPa'
$ws.Range("K21").Value = '13.5 MJ/m2'
$ws.Range("M21").Value = '20.7 °C 14:20 TU'
$ws.Range("O21").Value = '7.6 °C'

$ws.Range("E22").Value = '2026-02-23 14:49:03'
$ws.Range("K22").Value = '14.0 MJ/m2'
$ws.Range("M22").Value = '8.2 °C 12:17 TU'

$ws.Range("E23").Value = '2026-02-23 14:49:05'
$ws.Range("I23").Value = '0.2 mm'
$ws.Range("K23").Value = '13.1 MJ/m2'
$ws.Range("O23").Value = '3.2 °C'

$ws.Range("E24").Value = '2026-02-23 14:49:08'
$ws.Range("H24").Value = "'89%"
$ws.Range("J24").Value = '1027.0 hPa'
$ws.Range("K24").Value = '13.2 MJ/m2'
$ws.Range("L24").Value = '13.3 km/h - 217º 14:14 TU'
$ws.Range("M24").Value = '17.5 °C 14:28 TU'
$ws.Range("O24").Value = '6.7 °C'

$ws.Range("E25").Value = '2026-02-23 14:49:10'
$ws.Range("K25").Value = '14.3 MJ/m2'
$ws.Range("M25").Value = '10.3 °C 14:28 TU'
$ws.Range("O25").Value = '5.7 °C'

$ws.Range("E26").Value = '2026-02-23 14:49:12'
$ws.Range("J26").Value = '1024.2 hPa'
$ws.Range("K26").Value = '13.3 MJ/m2'
$ws.Range("O26").Value = '10.0 °C'

$ws.Range("E27").Value = '2026-02-23 14:49:15'
$ws.Range("K27").Value = '14.0 MJ/m2'
$ws.Range("O27").Value = '5.7 °C'

$ws.Range("E28").Value = '2026-02-23 14:49:17'
$ws.Range("H28").Value = "'71%"
$ws.Range("J28").Value = '1025.7 hPa'
$ws.Range("K28").Value = '12.7 MJ/m2'
$ws.Range("L28").Value = '21.6 km/h - 25º 14:25 TU'
$ws.Range("O28").Value = '9.7 °C'

$ws.Range("E29").Value = '2026-02-23 14:49:19'
$ws.Range("H29").Value = "'81%"
$ws.Range("K29").Value = '13.2 MJ/m2'
$ws.Range("O29").Value = '9.9 °C'

$ws.Range("E30").Value = '2026-02-23 14:49:22'
$ws.Range("H30").Value = "'69%"
$ws.Range("J30").Value = '1025.0 hPa'
$ws.Range("K30").Value = '13.0 MJ/m2'
$ws.Range("L30").Value = '19.1 km/h - 178º 14:24 TU'
$ws.Range("O30").Value = '12.4 °C'

$ws.Range("E31").Value = '2026-02-23 14:49:24'
$ws.Range("J31").Value = '1024.3 hPa'
$ws.Range("K31").Value = '13.0 MJ/m2'

$ws.Range("E32").Value = '2026-02-23 14:49:27'
$ws.Range("H32").Value = "'71%"
$ws.Range("K32").Value = '13.0 MJ/m2'
$ws.Range("O32").Value = '7.3 °C'

$ws.Range("E33").Value = '2026-02-23 14:49:29'
$ws.Range("H33").Value = "'51%"
$ws.Range("J33").Value = '1026.8 hPa'
$ws.Range("K33").Value = '13.5 MJ/m2'
$ws.Range("M33").Value = '18.9 °C 14:12 TU'
$ws.Range("O33").Value = '7.0 °C'

$ws.Range("E34").Value = '2026-02-23 14:49:31'
$ws.Range("K34").Value = '12.9 MJ/m2'
$ws.Range("O34").Value = '4.4 °C'

$ws.Range("E35").Value = '2026-02-23 14:49:34'
$ws.Range("J35").Value = '1025.5 hPa'
$ws.Range("K35").Value = '13.8 MJ/m2'
$ws.Range("M35").Value = '18.3 °C 14:06 TU'
$ws.Range("O35").Value = '12.3 °C'

$ws.Range("E36").Value = '2026-02-23 14:49:36'
$ws.Range("H36").Value = "'78%"
$ws.Range("K36").Value = '12.9 MJ/m2'
$ws.Range("O36").Value = '11.7 °C'

$ws.Range("E37").Value = '2026-02-23 14:49:39'
$ws.Range("H37").Value = "'65%"
$ws.Range("J37").Value = '1027.4 hPa'
$ws.Range("M37").Value = '17.7 °C 14:00 TU'
$ws.Range("O37").Value = '8.3 °C'

$ws.Range("E38").Value = '2026-02-23 14:49:41'
$ws.Range("H38").Value = "'63%"
$ws.Range("K38").Value = '13.4 MJ/m2'
$ws.Range("O38").Value = '11.2 °C'

$ws.Range("E39").Value = '2026-02-23 14:49:43'
$ws.Range("K39").Value = '14.7 MJ/m2'

$ws.Range("E40").Value = '2026-02-23 14:49:46'
$ws.Range("H40").Value = "'68%"
$ws.Range("J40").Value = '1027.4 hPa'
$ws.Range("O40").Value = '7.3 °C'

$ws.Range("E41").Value = '2026-02-23 14:49:48'
$ws.Range("K41").Value = '13.1 MJ/m2'
$ws.Range("O41").Value = '11.6 °C'

$ws.Range("E42").Value = '2026-02-23 14:49:50'
$ws.Range("H42").Value = "'81%"
$ws.Range("O42").Value = '10.6 °C'

$ws.Range("E43").Value = '2026-02-23 14:49:52'
$ws.Range("H43").Value = "'76%"
$ws.Range("K43").Value = '13.0 MJ/m2'
$ws.Range("L43").Value = '10.1 km/h - 116º 14:12 TU'
$ws.Range("M43").Value = '21.2 °C 14:00 TU'
$ws.Range("O43").Value = '8.4 °C'

$ws.Range("E44").Value = '2026-02-23 14:49:55'
$ws.Range("K44").Value = '13.4 MJ/m2'
$ws.Range("M44").Value = '5.6 °C 14:08 TU'
$ws.Range("O44").Value = '3.3 °C'

$ws.Range("E45").Value = '2026-02-23 14:49:57'
$ws.Range("J45").Value = '1027.6 hPa'
$ws.Range("K45").Value = '12.7 MJ/m2'
$ws.Range("L45").Value = '18.0 km/h - 206º 14:12 TU'
$ws.Range("O45").Value = '8.3 °C'

$ws.Range("E46").Value = '2026-02-23 14:50:00'
$ws.Range("H46").Value = "'81%"
$ws.Range("J46").Value = '1026.9 hPa'
$ws.Range("K46").Value = '12.9 MJ/m2'
$ws.Range("M46").Value = '22.2 °C 14:28 TU'
$ws.Range("O46").Value = '7.7 °C'
